# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 37
$ws1.Range("F4").Value = 223
$ws1.Range("F5").Value = 2688
$ws1.Range("F6").Value = 1900
$ws1.Range("F7").Value = 366
$ws1.Range("F8").Value = 115
$ws1.Range("F9").Value = 948

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 37
$ws4.Range("F4").Value = 223
$ws4.Range("F5").Value = 2688
$ws4.Range("F6").Value = 1900
$ws4.Range("F7").Value = 366
$ws4.Range("F9").Value = 115
$ws4.Range("F10").Value = 948
